# Auto-generated-by-analysis edit script: apply new TPM data to Angpt1-Tek.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-obsolete last data row (old row 17); rows 2-16 get fully
# overwritten below with the recomputed TPM values, and removing row 17
# shifts the sheet dimension down from A1:T17 to A1:T16.
$ws.Rows.Item(17).Delete() | Out-Null

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Angpt1"
$ws.Cells.Item(2,3).Value = "Tek"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = [double]"2"
$ws.Cells.Item(2,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(2,7).Value = [double]"0.2366953333333333"
$ws.Cells.Item(2,8).Value = [double]"0.710086"
$ws.Cells.Item(2,9).Value = [double]"0.01942010494447644"
$ws.Cells.Item(2,10).Value = [double]"0.01942010494447644"
$ws.Cells.Item(2,11).Value = [double]"3"
$ws.Cells.Item(2,12).Value = [double]"1"
$ws.Cells.Item(2,13).Value = [double]"58.51417433333334"
$ws.Cells.Item(2,14).Value = [double]"175.542523"
$ws.Cells.Item(2,15).Value = [double]"0.9208013000516164"
$ws.Cells.Item(2,16).Value = [double]"0.9208013000516164"
$ws.Cells.Item(2,17).Value = [double]"13.85003199855311"
$ws.Cells.Item(2,18).Value = [double]"124.650287986978"
$ws.Cells.Item(2,19).Value = [double]"0.01788205788001273"
$ws.Cells.Item(2,20).Value = [double]"0.01788205788001273"

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Angpt1"
$ws.Cells.Item(3,3).Value = "Tek"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = [double]"2"
$ws.Cells.Item(3,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(3,7).Value = [double]"0.2366953333333333"
$ws.Cells.Item(3,8).Value = [double]"0.710086"
$ws.Cells.Item(3,9).Value = [double]"0.01942010494447644"
$ws.Cells.Item(3,10).Value = [double]"0.01942010494447644"
$ws.Cells.Item(3,11).Value = [double]"3"
$ws.Cells.Item(3,12).Value = [double]"1"
$ws.Cells.Item(3,13).Value = [double]"4.208408333333334"
$ws.Cells.Item(3,14).Value = [double]"12.625225"
$ws.Cells.Item(3,15).Value = [double]"0.06622511397676659"
$ws.Cells.Item(3,16).Value = [double]"0.06622511397676657"
$ws.Cells.Item(3,17).Value = [double]"0.9961106132611113"
$ws.Cells.Item(3,18).Value = [double]"8.964995519349999"
$ws.Cells.Item(3,19).Value = [double]"0.001286098663388721"
$ws.Cells.Item(3,20).Value = [double]"0.001286098663388721"

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Angpt1"
$ws.Cells.Item(4,3).Value = "Tek"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = [double]"2"
$ws.Cells.Item(4,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(4,7).Value = [double]"0.2366953333333333"
$ws.Cells.Item(4,8).Value = [double]"0.710086"
$ws.Cells.Item(4,9).Value = [double]"0.01942010494447644"
$ws.Cells.Item(4,10).Value = [double]"0.01942010494447644"
$ws.Cells.Item(4,11).Value = [double]"2"
$ws.Cells.Item(4,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(4,13).Value = [double]"0.1957283333333334"
$ws.Cells.Item(4,14).Value = [double]"0.5871850000000001"
$ws.Cells.Item(4,15).Value = [double]"0.003080055488155473"
$ws.Cells.Item(4,16).Value = [double]"0.003080055488155473"
$ws.Cells.Item(4,17).Value = [double]"0.04632798310111112"
$ws.Cells.Item(4,18).Value = [double]"0.4169518479100001"
$ws.Cells.Item(4,19).Value = [double]"5.981500081478991E-05"
$ws.Cells.Item(4,20).Value = [double]"5.98150008147899E-05"

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Angpt1"
$ws.Cells.Item(5,3).Value = "Tek"
$ws.Cells.Item(5,4).Value = "Neutrophils"
$ws.Cells.Item(5,5).Value = [double]"2"
$ws.Cells.Item(5,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(5,7).Value = [double]"0.2366953333333333"
$ws.Cells.Item(5,8).Value = [double]"0.710086"
$ws.Cells.Item(5,9).Value = [double]"0.01942010494447644"
$ws.Cells.Item(5,10).Value = [double]"0.01942010494447644"
$ws.Cells.Item(5,11).Value = [double]"3"
$ws.Cells.Item(5,12).Value = [double]"1"
$ws.Cells.Item(5,13).Value = [double]"0.6056613333333333"
$ws.Cells.Item(5,14).Value = [double]"1.816984"
$ws.Cells.Item(5,15).Value = [double]"0.00953091707228673"
$ws.Cells.Item(5,16).Value = [double]"0.00953091707228673"
$ws.Cells.Item(5,17).Value = [double]"0.1433572111804444"
$ws.Cells.Item(5,18).Value = [double]"1.290214900624"
$ws.Cells.Item(5,19).Value = [double]"0.0001850914097609105"
$ws.Cells.Item(5,20).Value = [double]"0.0001850914097609105"

# Row 6
$ws.Cells.Item(6,1).Value = "ECs"
$ws.Cells.Item(6,2).Value = "Angpt1"
$ws.Cells.Item(6,3).Value = "Tek"
$ws.Cells.Item(6,4).Value = "Resolving-Mac"
$ws.Cells.Item(6,5).Value = [double]"2"
$ws.Cells.Item(6,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(6,7).Value = [double]"0.2366953333333333"
$ws.Cells.Item(6,8).Value = [double]"0.710086"
$ws.Cells.Item(6,9).Value = [double]"0.01942010494447644"
$ws.Cells.Item(6,10).Value = [double]"0.01942010494447644"
$ws.Cells.Item(6,11).Value = [double]"1"
$ws.Cells.Item(6,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(6,13).Value = [double]"0.023043"
$ws.Cells.Item(6,14).Value = [double]"0.069129"
$ws.Cells.Item(6,15).Value = [double]"0.0003626134111748421"
$ws.Cells.Item(6,16).Value = [double]"0.0003626134111748421"
$ws.Cells.Item(6,17).Value = [double]"0.005454170566"
$ws.Cells.Item(6,18).Value = [double]"0.049087535094"
$ws.Cells.Item(6,19).Value = [double]"7.041990499290021E-06"
$ws.Cells.Item(6,20).Value = [double]"7.041990499290021E-06"

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Angpt1"
$ws.Cells.Item(7,3).Value = "Tek"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = [double]"3"
$ws.Cells.Item(7,6).Value = [double]"1"
$ws.Cells.Item(7,7).Value = [double]"11.563232"
$ws.Cells.Item(7,8).Value = [double]"34.689696"
$ws.Cells.Item(7,9).Value = [double]"0.9487266849536318"
$ws.Cells.Item(7,10).Value = [double]"0.9487266849536319"
$ws.Cells.Item(7,11).Value = [double]"3"
$ws.Cells.Item(7,12).Value = [double]"1"
$ws.Cells.Item(7,13).Value = [double]"58.51417433333334"
$ws.Cells.Item(7,14).Value = [double]"175.542523"
$ws.Cells.Item(7,15).Value = [double]"0.9208013000516164"
$ws.Cells.Item(7,16).Value = [double]"0.9208013000516164"
$ws.Cells.Item(7,17).Value = [double]"676.6129731047787"
$ws.Cells.Item(7,18).Value = [double]"6089.516757943008"
$ws.Cells.Item(7,19).Value = [double]"0.8735887648989644"
$ws.Cells.Item(7,20).Value = [double]"0.8735887648989645"

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Angpt1"
$ws.Cells.Item(8,3).Value = "Tek"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = [double]"3"
$ws.Cells.Item(8,6).Value = [double]"1"
$ws.Cells.Item(8,7).Value = [double]"11.563232"
$ws.Cells.Item(8,8).Value = [double]"34.689696"
$ws.Cells.Item(8,9).Value = [double]"0.9487266849536318"
$ws.Cells.Item(8,10).Value = [double]"0.9487266849536319"
$ws.Cells.Item(8,11).Value = [double]"3"
$ws.Cells.Item(8,12).Value = [double]"1"
$ws.Cells.Item(8,13).Value = [double]"4.208408333333334"
$ws.Cells.Item(8,14).Value = [double]"12.625225"
$ws.Cells.Item(8,15).Value = [double]"0.06622511397676659"
$ws.Cells.Item(8,16).Value = [double]"0.06622511397676657"
$ws.Cells.Item(8,17).Value = [double]"48.66280190906667"
$ws.Cells.Item(8,18).Value = [double]"437.9652171816"
$ws.Cells.Item(8,19).Value = [double]"0.06282953284385419"
$ws.Cells.Item(8,20).Value = [double]"0.06282953284385419"

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Angpt1"
$ws.Cells.Item(9,3).Value = "Tek"
$ws.Cells.Item(9,4).Value = "MuSCs"
$ws.Cells.Item(9,5).Value = [double]"3"
$ws.Cells.Item(9,6).Value = [double]"1"
$ws.Cells.Item(9,7).Value = [double]"11.563232"
$ws.Cells.Item(9,8).Value = [double]"34.689696"
$ws.Cells.Item(9,9).Value = [double]"0.9487266849536318"
$ws.Cells.Item(9,10).Value = [double]"0.9487266849536319"
$ws.Cells.Item(9,11).Value = [double]"2"
$ws.Cells.Item(9,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(9,13).Value = [double]"0.1957283333333334"
$ws.Cells.Item(9,14).Value = [double]"0.5871850000000001"
$ws.Cells.Item(9,15).Value = [double]"0.003080055488155473"
$ws.Cells.Item(9,16).Value = [double]"0.003080055488155473"
$ws.Cells.Item(9,17).Value = [double]"2.263252127306667"
$ws.Cells.Item(9,18).Value = [double]"20.36926914576"
$ws.Cells.Item(9,19).Value = [double]"0.002922130832750983"
$ws.Cells.Item(9,20).Value = [double]"0.002922130832750982"

# Row 10
$ws.Cells.Item(10,1).Value = "FAPs"
$ws.Cells.Item(10,2).Value = "Angpt1"
$ws.Cells.Item(10,3).Value = "Tek"
$ws.Cells.Item(10,4).Value = "Neutrophils"
$ws.Cells.Item(10,5).Value = [double]"3"
$ws.Cells.Item(10,6).Value = [double]"1"
$ws.Cells.Item(10,7).Value = [double]"11.563232"
$ws.Cells.Item(10,8).Value = [double]"34.689696"
$ws.Cells.Item(10,9).Value = [double]"0.9487266849536318"
$ws.Cells.Item(10,10).Value = [double]"0.9487266849536319"
$ws.Cells.Item(10,11).Value = [double]"3"
$ws.Cells.Item(10,12).Value = [double]"1"
$ws.Cells.Item(10,13).Value = [double]"0.6056613333333333"
$ws.Cells.Item(10,14).Value = [double]"1.816984"
$ws.Cells.Item(10,15).Value = [double]"0.00953091707228673"
$ws.Cells.Item(10,16).Value = [double]"0.00953091707228673"
$ws.Cells.Item(10,17).Value = [double]"7.003402510762665"
$ws.Cells.Item(10,18).Value = [double]"63.03062259686399"
$ws.Cells.Item(10,19).Value = [double]"0.009042235358558564"
$ws.Cells.Item(10,20).Value = [double]"0.009042235358558564"

# Row 11
$ws.Cells.Item(11,1).Value = "FAPs"
$ws.Cells.Item(11,2).Value = "Angpt1"
$ws.Cells.Item(11,3).Value = "Tek"
$ws.Cells.Item(11,4).Value = "Resolving-Mac"
$ws.Cells.Item(11,5).Value = [double]"3"
$ws.Cells.Item(11,6).Value = [double]"1"
$ws.Cells.Item(11,7).Value = [double]"11.563232"
$ws.Cells.Item(11,8).Value = [double]"34.689696"
$ws.Cells.Item(11,9).Value = [double]"0.9487266849536318"
$ws.Cells.Item(11,10).Value = [double]"0.9487266849536319"
$ws.Cells.Item(11,11).Value = [double]"1"
$ws.Cells.Item(11,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(11,13).Value = [double]"0.023043"
$ws.Cells.Item(11,14).Value = [double]"0.069129"
$ws.Cells.Item(11,15).Value = [double]"0.0003626134111748421"
$ws.Cells.Item(11,16).Value = [double]"0.0003626134111748421"
$ws.Cells.Item(11,17).Value = [double]"0.266451554976"
$ws.Cells.Item(11,18).Value = [double]"2.398063994784"
$ws.Cells.Item(11,19).Value = [double]"0.0003440210195036362"
$ws.Cells.Item(11,20).Value = [double]"0.0003440210195036362"

# Row 12
$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(12,2).Value = "Angpt1"
$ws.Cells.Item(12,3).Value = "Tek"
$ws.Cells.Item(12,4).Value = "ECs"
$ws.Cells.Item(12,5).Value = [double]"3"
$ws.Cells.Item(12,6).Value = [double]"1"
$ws.Cells.Item(12,7).Value = [double]"0.3882320000000001"
$ws.Cells.Item(12,8).Value = [double]"1.164696"
$ws.Cells.Item(12,9).Value = [double]"0.03185321010189179"
$ws.Cells.Item(12,10).Value = [double]"0.03185321010189179"
$ws.Cells.Item(12,11).Value = [double]"3"
$ws.Cells.Item(12,12).Value = [double]"1"
$ws.Cells.Item(12,13).Value = [double]"58.51417433333334"
$ws.Cells.Item(12,14).Value = [double]"175.542523"
$ws.Cells.Item(12,15).Value = [double]"0.9208013000516164"
$ws.Cells.Item(12,16).Value = [double]"0.9208013000516164"
$ws.Cells.Item(12,17).Value = [double]"22.71707492977867"
$ws.Cells.Item(12,18).Value = [double]"204.4536743680081"
$ws.Cells.Item(12,19).Value = [double]"0.02933047727263925"
$ws.Cells.Item(12,20).Value = [double]"0.02933047727263925"

# Row 13
$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(13,2).Value = "Angpt1"
$ws.Cells.Item(13,3).Value = "Tek"
$ws.Cells.Item(13,4).Value = "FAPs"
$ws.Cells.Item(13,5).Value = [double]"3"
$ws.Cells.Item(13,6).Value = [double]"1"
$ws.Cells.Item(13,7).Value = [double]"0.3882320000000001"
$ws.Cells.Item(13,8).Value = [double]"1.164696"
$ws.Cells.Item(13,9).Value = [double]"0.03185321010189179"
$ws.Cells.Item(13,10).Value = [double]"0.03185321010189179"
$ws.Cells.Item(13,11).Value = [double]"3"
$ws.Cells.Item(13,12).Value = [double]"1"
$ws.Cells.Item(13,13).Value = [double]"4.208408333333334"
$ws.Cells.Item(13,14).Value = [double]"12.625225"
$ws.Cells.Item(13,15).Value = [double]"0.06622511397676659"
$ws.Cells.Item(13,16).Value = [double]"0.06622511397676657"
$ws.Cells.Item(13,17).Value = [double]"1.633838784066667"
$ws.Cells.Item(13,18).Value = [double]"14.7045490566"
$ws.Cells.Item(13,19).Value = [double]"0.002109482469523677"
$ws.Cells.Item(13,20).Value = [double]"0.002109482469523676"

# Row 14
$ws.Cells.Item(14,1).Value = "MuSCs"
$ws.Cells.Item(14,2).Value = "Angpt1"
$ws.Cells.Item(14,3).Value = "Tek"
$ws.Cells.Item(14,4).Value = "MuSCs"
$ws.Cells.Item(14,5).Value = [double]"3"
$ws.Cells.Item(14,6).Value = [double]"1"
$ws.Cells.Item(14,7).Value = [double]"0.3882320000000001"
$ws.Cells.Item(14,8).Value = [double]"1.164696"
$ws.Cells.Item(14,9).Value = [double]"0.03185321010189179"
$ws.Cells.Item(14,10).Value = [double]"0.03185321010189179"
$ws.Cells.Item(14,11).Value = [double]"2"
$ws.Cells.Item(14,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(14,13).Value = [double]"0.1957283333333334"
$ws.Cells.Item(14,14).Value = [double]"0.5871850000000001"
$ws.Cells.Item(14,15).Value = [double]"0.003080055488155473"
$ws.Cells.Item(14,16).Value = [double]"0.003080055488155473"
$ws.Cells.Item(14,17).Value = [double]"0.07598800230666669"
$ws.Cells.Item(14,18).Value = [double]"0.6838920207600002"
$ws.Cells.Item(14,19).Value = [double]"9.810965458970119E-05"
$ws.Cells.Item(14,20).Value = [double]"9.810965458970117E-05"

# Row 15
$ws.Cells.Item(15,1).Value = "MuSCs"
$ws.Cells.Item(15,2).Value = "Angpt1"
$ws.Cells.Item(15,3).Value = "Tek"
$ws.Cells.Item(15,4).Value = "Neutrophils"
$ws.Cells.Item(15,5).Value = [double]"3"
$ws.Cells.Item(15,6).Value = [double]"1"
$ws.Cells.Item(15,7).Value = [double]"0.3882320000000001"
$ws.Cells.Item(15,8).Value = [double]"1.164696"
$ws.Cells.Item(15,9).Value = [double]"0.03185321010189179"
$ws.Cells.Item(15,10).Value = [double]"0.03185321010189179"
$ws.Cells.Item(15,11).Value = [double]"3"
$ws.Cells.Item(15,12).Value = [double]"1"
$ws.Cells.Item(15,13).Value = [double]"0.6056613333333333"
$ws.Cells.Item(15,14).Value = [double]"1.816984"
$ws.Cells.Item(15,15).Value = [double]"0.00953091707228673"
$ws.Cells.Item(15,16).Value = [double]"0.00953091707228673"
$ws.Cells.Item(15,17).Value = [double]"0.2351371107626667"
$ws.Cells.Item(15,18).Value = [double]"2.116233996864"
$ws.Cells.Item(15,19).Value = [double]"0.0003035903039672567"
$ws.Cells.Item(15,20).Value = [double]"0.0003035903039672567"

# Row 16
$ws.Cells.Item(16,1).Value = "MuSCs"
$ws.Cells.Item(16,2).Value = "Angpt1"
$ws.Cells.Item(16,3).Value = "Tek"
$ws.Cells.Item(16,4).Value = "Resolving-Mac"
$ws.Cells.Item(16,5).Value = [double]"3"
$ws.Cells.Item(16,6).Value = [double]"1"
$ws.Cells.Item(16,7).Value = [double]"0.3882320000000001"
$ws.Cells.Item(16,8).Value = [double]"1.164696"
$ws.Cells.Item(16,9).Value = [double]"0.03185321010189179"
$ws.Cells.Item(16,10).Value = [double]"0.03185321010189179"
$ws.Cells.Item(16,11).Value = [double]"1"
$ws.Cells.Item(16,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(16,13).Value = [double]"0.023043"
$ws.Cells.Item(16,14).Value = [double]"0.069129"
$ws.Cells.Item(16,15).Value = [double]"0.0003626134111748421"
$ws.Cells.Item(16,16).Value = [double]"0.0003626134111748421"
$ws.Cells.Item(16,17).Value = [double]"0.008946029976000001"
$ws.Cells.Item(16,18).Value = [double]"0.080514269784"
$ws.Cells.Item(16,19).Value = [double]"1.155040117191592E-05"
$ws.Cells.Item(16,20).Value = [double]"1.155040117191592E-05"

